$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Last": add the B:E data/formula block (rows 1-8) that mirrors the
# layout already present on sheet "Main", plus keep the pre-existing A3
# "ignored error" cell intact.
# ---------------------------------------------------------------------------
$wsLast = $wb.Worksheets.Item("Last")

# Row 1 - plain (non-shared) formulas, same as the equivalent row on "Main".
$wsLast.Range("B1").Value = 1
$wsLast.Range("C1").Formula = "=B1+1"
$wsLast.Range("D1").Formula = "=C1+2"
$wsLast.Range("E1").Formula = "=D1+1"

# Rows 2-8 - column B values plus a horizontal C:E fill per row (creates a
# shared formula group per row, e.g. C2:E2, C3:E3, ...).
for ($r = 2; $r -le 8; $r++) {
    $wsLast.Range("B$r").Value = $r
    $wsLast.Range("C$r`:E$r").Formula = "=B$r+1"
}

# Column D, rows 2-8 - a single vertical shared formula group that overwrites
# the D cells from the row fills above (written last, starting at D2, so it
# cleanly becomes its own group instead of merging into the C:E ones).
$wsLast.Range("D2:D8").Formula = "=C2+2"

# Selection left on this sheet once data entry is done.
$wsLast.Range("A9").Select()

# ---------------------------------------------------------------------------
# Sheet "Main": selection moves from A13 to C1; re-activate it afterwards so
# it remains the active/visible tab (as in the original workbook).
# ---------------------------------------------------------------------------
$wsMain = $wb.Worksheets.Item("Main")
$wsMain.Activate()
$wsMain.Range("C1").Select()
